# Apply the changes described by the diff:
#  1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
#  2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on "Monthly Trend" sheet
#  3. Add a new "PO Forecast" worksheet (after the existing sheets) with forecast data

$wb = $excel.ActiveWorkbook

# --- 1) Weekly Quantity sheet -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2) Monthly Trend sheet ---------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3) New PO Forecast sheet -------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row (values first, then copy the existing header formatting over)
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$dates = @(45578.99999999999, 45599.99999999999, 45606.99999999999, 45613.99999999999, 45620.99999999999, 45627.99999999999, 45634.99999999999, 45641.99999999999, 45648.99999999999, 45655.99999999999)
$forecasts = @(240, 108, 64, 20, 0, 0, 0, 0, 0, 0)
$lowers = @(239.998952708904, 107.9989527379748, 63.9989525355079, 19.9989518827445, -24.00104917624263, -68.00105034958443, -112.001051693901, -156.0010533082734, -200.0010550218065, -244.0010569051364)
$uppers = @(239.9989527214002, 107.9989527504708, 63.9989529666653, 19.99895360798815, -24.0010454223651, -68.00104422103568, -112.0010427351548, -156.0010412979094, -200.0010394705959, -244.0010376102598)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value = $dates[$i]
    $wsForecast.Cells.Item($row, 2).Value = $forecasts[$i]
    $wsForecast.Cells.Item($row, 3).Value = $lowers[$i]
    $wsForecast.Cells.Item($row, 4).Value = $uppers[$i]
}

# Copy the date-format style from the source sheets onto column A of the data rows
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)
